$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.887.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "'1.631.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'214.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'0.5109"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.2547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "'0.06328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("D10").Value = "'19.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'4.259"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'1.635.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'0.5402"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'0.0₅7703"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "'63.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "'25.879.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'4.407"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'194.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("D21").Value = "'9.881"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'6.010"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'1.858"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "
$ws.Range("D25").Value = "'140.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").Value = "'0.1187"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.69%  "
$ws.Range("D27").Value = "'6.802"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").Value = "'15.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "'1.234"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "'0.04887"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Value = "'3.232"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'1.524"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").Value = "'2.363"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'0.8852"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'2.573"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "'1.132.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("D38").Value = "'0.5375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.16%  "
$ws.Range("D39").Value = "'0.01543"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "'0.8105"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "'0.0₈124"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.35%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'5.442"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'98.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "'1.768.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'0.4524"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'54.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Value = "'0.05048"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "'1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.36%  "
